$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers: insert "Reason" before "Food Items" (shift Food Items to J1)
$ws.Range("I1").Value = "Reason"
$ws.Range("J1").Value = "Food Items"

# Data rows (Bill No, Date, Cashier, KOT, Price, SGST, CGST, Tax, Reason, Food Items)
$rows = @(
    @(412, 45724.22928240741, "Ajay Francis Anchan", 4, 220, 2.5, 2.5, 5, "fdg", "Chicken Burger (x1), Vanilla Shake (x4)"),
    @(391, 45721.22928240741, "Ajay Francis Anchan", 3, 300, 4.5, 4.5, 9, "fff", "Chicken Burger (x5)"),
    @(389, 45721.22928240741, "Ajay Francis Anchan", 1, 380, 4.35, 4.35, 10.5, "er", "Chicken Cheese Burger (x1), Vanilla Shake (x1), Butterscotch Lassi (x1), Chicken Wrap (x2), Oreo Shake (x2)"),
    @(388, 45721.22928240741, "Ajay Francis Anchan", 2, 60, 0.9, 0.9, 1.8, "drg", "Mango Lassi (x1)"),
    @(387, 45721.22928240741, "Ajay Francis Anchan", 1, 60, 0, 0, 1.8, "no", "Chicken Cheese Burger (x1)"),
    @(386, 45721.22928240741, "Ajay Francis Anchan", 3, 120, 1.2, 1.2, 2.4, "sdf", "Strawberry Lassi (x3)"),
    @(385, 45721.22928240741, "Ajay Francis Anchan", 2, 510, 7.58, 7.58, 15.15, "hate this", "Vanilla Shake (x3), Chicken Wrap (x3), Veg Cheese Pops (x3)"),
    @(384, 45721.22928240741, "Ajay Francis Anchan", 1, 200, 0.8, 0.8, 5.2, "sf", "Chicken Cheese Burger (x2), Vanilla Shake (x2)"),
    @(364, 45719.22928240741, "Ajay Francis Anchan", 7, 110, 1.4, 1.4, 2.8, "mjkf", "Mango Lassi (x1), Banana Shake (x1)"),
    @(342, 45718.22928240741, "Ajay Francis Anchan", 93, 60, 0.9, 0.9, 1.8, "jjjk", "Chicken Burger (x1)"),
    @(252, 45718.22928240741, "Ajay Francis Anchan", 3, 280, 4, 4, 8, "rsfgff", "Mango Lassi (x1), Chicken Burger (x1), Strawberry Lassi (x1), Butterscotch Lassi (x1), Chicken Wrap (x1)"),
    @(222, 45709.22928240741, "Ajay Francis Anchan", 3, 160, 2.2, 2.2, 4.4, "I hate this", "Strawberry Lassi (x1), Mango Lassi (x1), Butterscotch Lassi (x1)"),
    @(220, 45709.22928240741, "Ajay Francis Anchan", 1, 120, 1.8, 1.8, 3.6, "I love this", "Chicken Burger (x2)")
)

# Copy the date format from the existing B2 cell so new date cells share style index 1
$ws.Range("B2").Copy()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $r++
}

# Apply the date number format to the whole B2:B14 range (PasteSpecial formats only,
# so the values just written above are preserved)
$ws.Range("B2:B14").PasteSpecial(-4122)
